# Applies the changes described by the commit:
# "Terminado filtros tabal consulta. Cargar Funcionario con tabal Auxiliar.
#  Base MySQL ausentismos_v2.sql."
#
# Concretely: widen columns C:D (NOMBRE / FEC_INGRESO) and append one new
# employee record (row 22) to the "Funcinarios Unicauca 2022" sheet, then
# leave the selection on the last entered cell (J22), matching the target
# workbook produced by the author's Excel session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen columns C (NOMBRE) and D (FEC_INGRESO) ---------------------
$ws.Columns("C:D").ColumnWidth = 18

# --- Append the new record in row 22 -----------------------------------
# Start from a duplicate of row 21 so the new row inherits the same cell
# formatting (in particular the date style used by D21/O21), then drop the
# two incidental empty cells (Q/S) that row 21 itself does not carry and
# overwrite the fields that actually differ for the new employee.
$ws.Range("A21:W21").Copy($ws.Range("A22:W22")) | Out-Null
$ws.Range("Q22").ClearContents() | Out-Null
$ws.Range("S22").ClearContents() | Out-Null

$ws.Range("A22").Value = 3444            # CODIGO
$ws.Range("B22").Value = 2020055         # CEDULA
$ws.Range("C22").Value = "MESSI g"       # NOMBRE
$ws.Range("F22").Value = "MAS"           # GENERO
$ws.Range("J22").Value = 10              # C_COSTO
$ws.Range("K22").Value = "DEPTO D"       # DEPARTAMENTO
$ws.Range("W22").Value = "ACTIVO"        # ESTADO

# --- Match the author's final selection --------------------------------
$ws.Range("J22").Select() | Out-Null
